# Staging.ReportingPeriod.xlsx - add FirstCycleDate, LastCycleDate and
# YearName columns to the header row, keeping all column headers in
# alphabetical order (A2:J2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "BusinessKey"
$ws.Range("B2").Value = "EndDateID"
$ws.Range("C2").Value = "FirstCycleDate"
$ws.Range("D2").Value = "ID"
$ws.Range("E2").Value = "LastCycleDate"
$ws.Range("F2").Value = "ReportingPeriod"
$ws.Range("G2").Value = "StartDateID"
$ws.Range("H2").Value = "Summary"
$ws.Range("I2").Value = "YearName"
$ws.Range("J2").Value = "YearNumber"
